# Update the date and all the two-digit multiplication problems in the
# document to match the new "output generated at c8c62b6" content.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-05-16 Friday"; New = "2025-05-17 Saturday" },
    @{ Old = "64×53=";            New = "87×33=" },
    @{ Old = "62×47=";            New = "60×73=" },
    @{ Old = "34×43=";            New = "51×84=" },
    @{ Old = "36×87=";            New = "84×27=" },
    @{ Old = "89×27=";            New = "67×44=" },
    @{ Old = "25×94=";            New = "36×51=" },
    @{ Old = "80×97=";            New = "13×37=" },
    @{ Old = "75×79=";            New = "41×65=" },
    @{ Old = "79×48=";            New = "81×38=" },
    @{ Old = "40×94=";            New = "51×17=" },
    @{ Old = "24×37=";            New = "31×14=" },
    @{ Old = "33×49=";            New = "22×44=" },
    @{ Old = "72×60=";            New = "43×16=" },
    @{ Old = "88×95=";            New = "84×84=" },
    @{ Old = "66×61=";            New = "67×75=" },
    @{ Old = "68×59=";            New = "50×13=" },
    @{ Old = "17×46=";            New = "87×79=" },
    @{ Old = "87×39=";            New = "33×80=" },
    @{ Old = "67×56=";            New = "91×21=" },
    @{ Old = "57×31=";            New = "64×11=" },
    @{ Old = "67×68=";            New = "24×22=" },
    @{ Old = "22×77=";            New = "43×48=" },
    @{ Old = "42×64=";            New = "60×64=" },
    @{ Old = "69×79=";            New = "53×58=" },
    @{ Old = "18×44=";            New = "65×50=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $r.New, 2) | Out-Null
}

$d.Save()
